# Applies:
#  1. Footer "date" field text 8/10/2021 -> 8/12/2021 on the Slide Master
#     and on every slide layout (the datetimeFigureOut field placeholder).
#  2. Re-positions three shapes on slide 3 ("Rectangle 3", "Rectangle 4",
#     "TextBox 5") to their new x/y offsets.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Update the "Date Placeholder" (ppPlaceholderDate = 16) text on the
#    slide master and on every custom layout.
# ---------------------------------------------------------------------
$oldDate = "8/10/2021"
$newDate = "8/12/2021"
$ppPlaceholderDate = 16

$master = $p.SlideMaster

for ($i = 1; $i -le $master.Shapes.Placeholders.Count; $i++) {
    $shp = $master.Shapes.Placeholders.Item($i)
    if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Placeholders.Count; $i++) {
        $shp = $layout.Shapes.Placeholders.Item($i)
        if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2. Move three shapes on slide 3 ("JavaScript Execution Environment").
#    Left/Top are expressed in points (EMU / 12700); the values below
#    were chosen so the stored EMU offsets match the target exactly.
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)

$rect3 = $slide3.Shapes.Item("Rectangle 3")
$rect3.Left = 196.9273308346457
$rect3.Top = 405.26820397637795

$rect4 = $slide3.Shapes.Item("Rectangle 4")
$rect4.Left = 196.92723909448821
$rect4.Top = 368.0566869133858

$textbox5 = $slide3.Shapes.Item("TextBox 5")
$textbox5.Left = 547.9834905669292
$textbox5.Top = 103.88157280314961
